$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with new column names
$ws.Range("B1").Value = "Chofer"
$ws.Range("C1").Value = "En_camino_hoy"
$ws.Range("D1").Value = "Entregados"
$ws.Range("E1").Value = "Nadie_en_domicilio"
$ws.Range("F1").Value = "No_visitado"
$ws.Range("G1").Value = "Diferencia_cargados_Entregados"
$ws.Range("H1").Value = "Horario_salida"
$ws.Range("I1").Value = "Horario_fin"

# New column J1, matching the existing header style
$ws.Range("J1").Value = "efectividad"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Remove the old data row 2 entirely
$ws.Rows.Item(2).Delete()
